# Adds a new "2021" column (R) to the 7.2.1.1 table:
#   R4 = 2021            (year header, same look as the other year headers)
#   R5 = 31.8             (growth-rate row)
#   R6 = 12957.1           (absolute value row)
# and moves the active selection onto the new column, matching the
# published diff (dimension A1:Q7 -> A1:R7, selection P7 -> R4:R6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R4: new year header cell -> copy formatting from the previous year header (Q4) ---
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value = 2021

# --- R5: growth-rate figure -> copy formatting from a General-formatted sibling (N5) ---
$ws.Range("N5").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R5").Value = 31.8

# --- R6: absolute-value figure -> copy formatting from the previous year's cell (Q6) ---
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R6").Value = 12957.1

# Clear the marching-ants copy indicator now that the paste operations are done.
$ws.Application.CutCopyMode = $false

# Match the committed selection state: active cell R4, selected range R4:R6.
$ws.Range("R4:R6").Select() | Out-Null
